# Apply the "ph_integrated_template" revision:
#  - rename the "% non-trauma deaths" label (instructions!C2) to
#    "% non-trauma deaths*" (matching the footnoted label already used
#    elsewhere in the workbook)
#  - clear the now-unused duplicate threshold cells in columns M, O and P
#    of the instructions sheet (rows 3-7) - the ranges D:E / H:J / L / Q
#    already carry these thresholds, so M/O/P are blanked out
#  - rename the shared "Number of Observation" label to
#    "Number of Observations" (Data!B2 and Cat!B2)
#  - give instructions row 3 an explicit (custom) row height
#  - restore the active selections on each sheet

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("instructions")
$wsData         = $wb.Worksheets.Item("Data")
$wsCat          = $wb.Worksheets.Item("Cat")

# --- instructions sheet -----------------------------------------------
$wsInstructions.Range("C2").Value = "% non-trauma deaths*"

foreach ($r in 3..7) {
    $wsInstructions.Range("M$r").Value = $null
    $wsInstructions.Range("O$r").Value = $null
    $wsInstructions.Range("P$r").Value = $null
}

$wsInstructions.Rows.Item(3).RowHeight = 30

# --- Data / Cat sheets ---------------------------------------------------
$wsData.Range("B2").Value = "Number of Observations"
$wsCat.Range("B2").Value  = "Number of Observations"

# --- restore selections (set last = left as the active tab) -------------
$null = $wsData.Range("H8").Select()
$null = $wsCat.Range("J8").Select()
$null = $wsInstructions.Range("H13").Select()
